$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Delete the "tại HỆ THỐNG" block rows (rows 4 through 10)
$ws.Rows("4:10").Delete()

# Update Ngày công and Phụ cấp values
$ws.Range("B2").Value = 8
$ws.Range("B3").Value = 280000

# Delete the "Tổng lương tại HỆ THỐNG" row (now row 28 after previous shift)
$ws.Rows("28:28").Delete()

# Update the total salary figures
$ws.Range("B28").Value = -720000
$ws.Range("B31").Value = -720000
